$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.742.88"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.852.28"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.18"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4300"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3711"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07367"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8766"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.03"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.865.00"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.454"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.606"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06946"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "81.15"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009086"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9997"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.56"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.766.46"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.093"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.03"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +6.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.083.64"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.965"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.86%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.94"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.58"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.335"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.59"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -5.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.854"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08925"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7831"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.612"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.973"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.163"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +5.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.0000"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.112"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05432"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01961"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.841"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5220"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1684"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.768"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.664"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.40%  "
$ws.Range("E45").Value = "  +2.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "107.25"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4777"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06568"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.52%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.0000"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.666"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.846"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +6.03%  "
